# Insert a new data row at row 135 in the "Apio" sheet.
# This shifts the existing rows 135-179 down to 136-180, and we populate
# the newly inserted row 135 with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 135 (pushes rows 135..179 down to 136..180)
$ws.Rows.Item(135).Insert()

# Populate the new row 135 with its data
$ws.Range("A135").Value = 4
$ws.Range("B135").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C135").Value = "Los Lagos"
$ws.Range("D135").Value = 44559
$ws.Range("E135").Value = 10
$ws.Range("F135").Value = 100112017
$ws.Range("G135").Value = "Apio"
$ws.Range("H135").Value = "Americana (o)"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 20
$ws.Range("K135").Value = 12000
$ws.Range("L135").Value = 12500
$ws.Range("M135").Value = 12250
$ws.Range("N135").Value = "`$/docena de matas"
$ws.Range("O135").Value = "Región de Coquimbo"
$ws.Range("P135").Value = 2042
$ws.Range("Q135").Value = 6
$ws.Range("R135").Value = "Hortaliza"
